# Append 5 new daily rows (27-09-2021 .. 01-10-2021) to the EMBI spreads
# table on Sheet1, right after the existing last row (185, dated
# 24-09-2021). The sheet's used range grows from A1:P185 to A1:P190 and
# the "Serie" date strings are added to the shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 186; Date = "27-09-2021"; Values = @(322.3, 398, 199, 323, 84.90000000000001, 118.3, 11.7, 151.6, 517.1, 1631.1, 298.5, 293, 142, 359, 172) },
    @{ Row = 187; Date = "28-09-2021"; Values = @(323, 399, 199, 325, 84.2, 117.7, 7.3, 152.5, 518.6, 1597.4, 302.1, 298, 143, 361, 173) },
    @{ Row = 188; Date = "29-09-2021"; Values = @(321.5, 396, 203, 320, 86.59999999999999, 119.4, 12, 149.7, 506.1, 1589.1, 297.8, 299, 144, 356, 175) },
    @{ Row = 189; Date = "30-09-2021"; Values = @(324.4, 399, 204, 318, 87.09999999999999, 119.4, 13.5, 153.2, 507.1, 1606.6, 304.2, 301, 150, 360, 180) },
    @{ Row = 190; Date = "01-10-2021"; Values = @(329.9, 403, 207, 326, 88.2, 121.2, 15.8, 160.7, 516.8, 1610.6, 309.5, 302, 155, 363, 186) }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Column A holds the "Serie" date as plain text (e.g. "24-09-2021" in
    # the row right above), not a real Excel date. Most of the new dates
    # (27/28/29/30-09-2021) aren't valid calendar dates anyway so Excel
    # leaves them as text automatically, but "01-10-2021" *is* a valid
    # date (1 Oct 2021) and would silently get auto-converted to a date
    # serial on entry. Force Text first so it is stored the same way as
    # every other cell in the column, then restore the default "Normal"
    # style so no stray number-format survives on the cell.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $entry.Date
    $dateCell.Style = "Normal"

    for ($i = 0; $i -lt $entry.Values.Count; $i++) {
        $col = $i + 2   # values start at column B
        $ws.Cells.Item($r, $col).Value = $entry.Values[$i]
    }
}
